# ----------------------------------------------------------------------
# This workbook's single data table (A1:D49) is re-sorted: each year's
# 12 monthly rows go from calendar order (Jan..Dec) to a fiscal order
# (Oct,Nov,Dec,Jan..Sep), and new rows for 2022 and part of 2023 are
# appended at the end, growing the sheet from A1:D49 to A1:D68.
# ----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every existing data row (A2:D49), keyed by its date in column A,
# before any cell gets overwritten -- rows are reshuffled both up and down,
# so a straight in-place copy would clobber data that is still needed.
$byDate = @{}
for ($r = 2; $r -le 49; $r++) {
    $rowDate = $ws.Cells.Item($r, 1).Value2
    $byDate[$rowDate] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2
    )
}

# Target layout: for each destination row, either re-use the snapshotted
# values for that month ($null placeholders below) or supply the literal
# values for newly added 2022/2023 months.
$targetRows = @(
    @(2, "2018-10", $null, $null, $null),
    @(3, "2018-11", $null, $null, $null),
    @(4, "2018-12", $null, $null, $null),
    @(5, "2018-01", $null, $null, $null),
    @(6, "2018-02", $null, $null, $null),
    @(7, "2018-03", $null, $null, $null),
    @(8, "2018-04", $null, $null, $null),
    @(9, "2018-05", $null, $null, $null),
    @(10, "2018-06", $null, $null, $null),
    @(11, "2018-07", $null, $null, $null),
    @(12, "2018-08", $null, $null, $null),
    @(13, "2018-09", $null, $null, $null),
    @(14, "2019-10", $null, $null, $null),
    @(15, "2019-11", $null, $null, $null),
    @(16, "2019-12", $null, $null, $null),
    @(17, "2019-01", $null, $null, $null),
    @(18, "2019-02", $null, $null, $null),
    @(19, "2019-03", $null, $null, $null),
    @(20, "2019-04", $null, $null, $null),
    @(21, "2019-05", $null, $null, $null),
    @(22, "2019-06", $null, $null, $null),
    @(23, "2019-07", $null, $null, $null),
    @(24, "2019-08", $null, $null, $null),
    @(25, "2019-09", $null, $null, $null),
    @(26, "2020-10", $null, $null, $null),
    @(27, "2020-11", $null, $null, $null),
    @(28, "2020-12", $null, $null, $null),
    @(29, "2020-01", $null, $null, $null),
    @(30, "2020-02", $null, $null, $null),
    @(31, "2020-03", $null, $null, $null),
    @(32, "2020-04", $null, $null, $null),
    @(33, "2020-05", $null, $null, $null),
    @(34, "2020-06", $null, $null, $null),
    @(35, "2020-07", $null, $null, $null),
    @(36, "2020-08", $null, $null, $null),
    @(37, "2020-09", $null, $null, $null),
    @(38, "2021-10", $null, $null, $null),
    @(39, "2021-11", $null, $null, $null),
    @(40, "2021-12", $null, $null, $null),
    @(41, "2021-01", $null, $null, $null),
    @(42, "2021-02", $null, $null, $null),
    @(43, "2021-03", $null, $null, $null),
    @(44, "2021-04", $null, $null, $null),
    @(45, "2021-05", $null, $null, $null),
    @(46, "2021-06", $null, $null, $null),
    @(47, "2021-07", $null, $null, $null),
    @(48, "2021-08", $null, $null, $null),
    @(49, "2021-09", $null, $null, $null),
    @(50, "2022-10", 90.2, 80.59999999999999, 92.40000000000001),
    @(51, "2022-11", 91, 84.09999999999999, 88.7),
    @(52, "2022-12", 88.7, 93.40000000000001, 90.59999999999999),
    @(53, "2022-01", 88.59999999999999, 92.90000000000001, 119.5),
    @(54, "2022-02", 122.4, 92.09999999999999, 120.9),
    @(55, "2022-03", 121, 91.8, 123.9),
    @(56, "2022-04", 121.1, 93.59999999999999, 123.1),
    @(57, "2022-05", 113.1, 87.3, 117.9),
    @(58, "2022-06", 141, 81, 117.6),
    @(59, "2022-07", 86.8, 73.40000000000001, 109.5),
    @(60, "2022-08", 95.09999999999999, 69.5, 102),
    @(61, "2022-09", 128, 74.09999999999999, 96.3),
    @(62, "2023-01", 158.6, 95.7, 90),
    @(63, "2023-02", 140.8, 94.59999999999999, 89.5),
    @(64, "2023-03", 130.7, 93.09999999999999, 85.3),
    @(65, "2023-04", 128.6, 89.40000000000001, 83),
    @(66, "2023-05", 122.7, 85.5, 83.90000000000001),
    @(67, "2023-06", 77.5, 86.2, 82.2),
    @(68, "2023-07", 115, 93, 87.09999999999999)
)

foreach ($item in $targetRows) {
    $r = $item[0]
    $rowDate = $item[1]
    if ($byDate.ContainsKey($rowDate)) {
        $vals = $byDate[$rowDate]
        $b = $vals[0]
        $c = $vals[1]
        $d = $vals[2]
    } else {
        $b = $item[2]
        $c = $item[3]
        $d = $item[4]
    }
    $ws.Cells.Item($r, 1).Value = $rowDate
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}
